# Formulaire Bricolage.docx edit script
# Converts the "Connexion/Inscription" account form into the
# "bricolage" (JP INFO BRICO SERVICES) service-request form.

$d = $word.ActiveDocument

# --- 1. Title paragraph: "Formulaire Connexion/Inscription" -> "Formulaire bricolage"
#        font size 14pt/28half-pt -> 20pt/40half-pt (both sz and szCs)
$titlePara = $d.Paragraphs.Item(2)
$titlePara.Range.Font.Size = 20
$titlePara.Range.Font.SizeBi = 20
$titlePara.Range.Find.Execute("Connexion/Inscription", $false, $false, $false, $false, $false, $true, 1, $false, "bricolage", 2) | Out-Null

# --- 2. Remove the label text from the account-only fields, keep the (now
#        empty) list paragraphs in place: Pseudo / Mot de Passe / Nom / Prenom
$pseudoPara = $d.Paragraphs.Item(3)
$r = $d.Range($pseudoPara.Range.Start, $pseudoPara.Range.End - 1)
$r.Text = ""

$passPara = $d.Paragraphs.Item(4)
$r = $d.Range($passPara.Range.Start, $passPara.Range.End - 1)
$r.Text = ""

$nomPara = $d.Paragraphs.Item(5)
$r = $d.Range($nomPara.Range.Start, $nomPara.Range.End - 1)
$r.Text = ""

$prenomPara = $d.Paragraphs.Item(6)
$r = $d.Range($prenomPara.Range.Start, $prenomPara.Range.End - 1)
$r.Text = ""

# --- 3. Relabel the remaining fields:
#        "Adresse :"      -> "Detail(s) : "
#        "Code Postal :"  -> "Adresse :"
#        "Ville :"        -> "Code Postal :"
#        "Telephone :"    -> "Ville" (keep the existing " :" run untouched)
$detailPara = $d.Paragraphs.Item(7)
$r = $d.Range($detailPara.Range.Start, $detailPara.Range.End - 1)
$r.Text = "Détail(s) : "

$adressePara = $d.Paragraphs.Item(8)
$r = $d.Range($adressePara.Range.Start, $adressePara.Range.End - 1)
$r.Text = "Adresse`u{00A0}:"

$cpPara = $d.Paragraphs.Item(9)
$r = $d.Range($cpPara.Range.Start, $cpPara.Range.End - 1)
$r.Text = "Code Postal`u{00A0}:"

$villePara = $d.Paragraphs.Item(10)
$villePara.Range.Find.Execute("Téléphone", $false, $false, $false, $false, $false, $true, 1, $false, "Ville", 2) | Out-Null

# --- 4. Drop the "Email :" field entirely (whole paragraph removed)
$emailPara = $d.Paragraphs.Item(11)
$emailPara.Range.Delete() | Out-Null
